$d = $word.ActiveDocument

$replacements = @(
    @{old="872×6=5232"; new="709×6=4254"},
    @{old="496×6=2976"; new="556×6=3336"},
    @{old="811×5=4055"; new="733×5=3665"},
    @{old="368×6=2208"; new="258×3=774"},
    @{old="159×2=318";  new="464×2=928"},
    @{old="236×4=944";  new="559×7=3913"},
    @{old="421×5=2105"; new="234×3=702"},
    @{old="252×6=1512"; new="706×7=4942"},
    @{old="471×7=3297"; new="749×2=1498"},
    @{old="148×6=888";  new="695×5=3475"},
    @{old="234×2=468";  new="668×4=2672"},
    @{old="718×2=1436"; new="745×3=2235"},
    @{old="848×2=1696"; new="347×2=694"},
    @{old="163×3=489";  new="444×7=3108"},
    @{old="894×6=5364"; new="474×8=3792"},
    @{old="971×9=8739"; new="139×6=834"},
    @{old="591×7=4137"; new="922×8=7376"},
    @{old="502×9=4518"; new="816×7=5712"},
    @{old="780×6=4680"; new="139×9=1251"},
    @{old="209×8=1672"; new="359×4=1436"},
    @{old="579×7=4053"; new="593×9=5337"},
    @{old="175×6=1050"; new="572×9=5148"},
    @{old="499×7=3493"; new="597×9=5373"},
    @{old="492×2=984";  new="712×2=1424"},
    @{old="800×8=6400"; new="128×8=1024"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2) | Out-Null
}
